$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.607.22"
$ws.Range("E2").Value = "  -1.04%  "

$ws.Range("D3").Value = "2.369.24"
$ws.Range("E3").Value = "  -1.60%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").Value = "'503.32"
$ws.Range("E5").Value = "  -1.74%  "

$ws.Range("D6").Value = "'132.38"
$ws.Range("E6").Value = "  -0.61%  "

$ws.Range("E7").Value = "  +0.80%  "

$ws.Range("D8").Value = "'0.551"
$ws.Range("E8").Value = "  -1.47%  "

$ws.Range("D9").Value = "2.366.65"
$ws.Range("E9").Value = "  -3.24%  "

$ws.Range("D10").Value = "'0.0980"
$ws.Range("E10").Value = "  +0.53%  "

$ws.Range("D11").Value = "'0.150"
$ws.Range("E11").Value = "  +0.23%  "

$ws.Range("D12").Value = "'0.330"
$ws.Range("E12").Value = "  +2.37%  "

$ws.Range("D13").Value = "'4.63"
$ws.Range("E13").Value = "  -1.67%  "

$ws.Range("D14").Value = "2.804.06"
$ws.Range("E14").Value = "  -1.15%  "

$ws.Range("D15").Value = "56.557.04"
$ws.Range("E15").Value = "  -0.77%  "

$ws.Range("D16").Value = "'21.57"
$ws.Range("E16").Value = "  -1.76%  "

$ws.Range("D17").Value = "'0.0000132"
$ws.Range("E17").Value = "  -1.15%  "

$ws.Range("D18").Value = "2.401.89"
$ws.Range("E18").Value = "  +1.09%  "

$ws.Range("D19").Value = "'10.03"
$ws.Range("E19").Value = "  -2.68%  "

$ws.Range("D20").Value = "'308.32"
$ws.Range("E20").Value = "  -1.78%  "

$ws.Range("D21").Value = "'4.01"
$ws.Range("E21").Value = "  -2.83%  "

$ws.Range("D22").Value = "'6.18"
$ws.Range("E22").Value = "  -4.14%  "

$ws.Range("E23").Value = "  +0.35%  "

$ws.Range("D24").Value = "'65.01"
$ws.Range("E24").Value = "  -0.17%  "

$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  +0.67%  "

$ws.Range("D26").Value = "'0.150"
$ws.Range("E26").Value = "  -0.62%  "

$ws.Range("E27").Value = "  -3.29%  "

$ws.Range("D28").Value = "'7.27"
$ws.Range("E28").Value = "  -2.97%  "

$ws.Range("D29").Value = "'171.23"
$ws.Range("E29").Value = "  -0.95%  "

$ws.Range("D30").Value = "0.0₃0721"
$ws.Range("E30").Value = "  -2.55%  "

$ws.Range("D31").Value = "'1.64"
$ws.Range("E31").Value = "  -2.83%  "

$ws.Range("D32").Value = "'1.11"
$ws.Range("E32").Value = "  -3.51%  "

$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "  +0.24%  "

$ws.Range("D34").Value = "'5.77"
$ws.Range("E34").Value = "  -7.30%  "

$ws.Range("D35").Value = "'0.998"
$ws.Range("E35").Value = "  +0.52%  "

$ws.Range("D36").Value = "'17.72"
$ws.Range("E36").Value = "  -1.95%  "

$ws.Range("E37").Value = "  -3.24%  "

$ws.Range("D38").Value = "'3.81"
$ws.Range("E38").Value = "  -1.10%  "

$ws.Range("D39").Value = "'0.811"
$ws.Range("E39").Value = "  -1.12%  "

$ws.Range("D40").Value = "'36.06"
$ws.Range("E40").Value = "  +0.00%  "

$ws.Range("E41").Value = "  -2.21%  "

$ws.Range("D42").Value = "'130.30"
$ws.Range("E42").Value = "  -0.59%  "

$ws.Range("D43").Value = "'3.36"
$ws.Range("E43").Value = "  -1.91%  "

$ws.Range("D44").Value = "'4.80"
$ws.Range("E44").Value = "  -3.96%  "

$ws.Range("D45").Value = "'0.561"
$ws.Range("E45").Value = "  -1.67%  "

$ws.Range("D46").Value = "'0.0907"
$ws.Range("E46").Value = "  -0.54%  "

$ws.Range("D47").Value = "'246.19"
$ws.Range("E47").Value = "  -5.43%  "

$ws.Range("D48").Value = "'0.0482"
$ws.Range("E48").Value = "  -3.03%  "

$ws.Range("D49").Value = "'0.0208"
$ws.Range("E49").Value = "  -2.28%  "

$ws.Range("D50").Value = "'17.06"
$ws.Range("E50").Value = "  -1.81%  "

$ws.Range("E51").Value = "  -1.21%  "
